# Prefix each protocol sheet's own name onto every "step/command" label
# in column A (skipping the header row "Name" in row 1).
#
# e.g. on sheet "discount2", cell A2 "Step4 Takeaway" becomes
# "discount2 Step4 Takeaway".

$wb = $excel.ActiveWorkbook

# These are the "protocol" worksheets whose Column A command names must be
# prefixed with the worksheet's own name. (All sheets except the first six
# overview/summary sheets.)
$protocolSheetNames = @(
    "price1", "price2",
    "discount1", "discount2",
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol", "dickpic", "boosters"
)

foreach ($sheetName in $protocolSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

    # Row 1 is the header ("Name") and must stay untouched; data starts row 2.
    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value2
        if ($current -ne $null -and $current -ne "") {
            $newValue = "$sheetName $current"
            $cell.Value = $newValue
        }
    }
}
